# Regenerate the "K" column (column G) values in save_data to reflect the
# newly computed std/mean based statistic (s_vals), replacing the older
# Strike#-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K (column G) value
$newValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 3
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    19 = 0
    20 = 2
    21 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
